# Apply updated crypto market data (price + 1h volume%) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.717.03"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "3.032.57"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'511.46"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'140.11"
$ws.Range("E6").Value = "  +4.04%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("D9").Value = "'7.50"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("D12").Value = "3.550.18"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "'26.72"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("E15").Value = "  +10.41%  "
$ws.Range("D16").Value = "57.744.52"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("E17").Value = "  +9.35%  "
$ws.Range("D18").Value = "3.037.08"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'12.89"
$ws.Range("E19").Value = "  +4.64%  "
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").Value = "'332.78"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "'5.83"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "'64.59"
$ws.Range("E25").Value = "  +4.79%  "
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "0.0₃0936"
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("E29").Value = "  +7.35%  "
$ws.Range("D30").Value = "'7.49"
$ws.Range("E30").Value = "  +11.25%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.81"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.22"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").Value = "'20.82"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +7.07%  "
$ws.Range("D35").Value = "'154.95"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  +7.11%  "
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("D38").Value = "'24.84"
$ws.Range("E38").Value = "  +8.13%  "
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "3.068.24"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").Value = "'37.42"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "'3.87"
$ws.Range("E42").Value = "  +9.32%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "2.309.86"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").Value = "'0.657"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("D46").Value = "'1.43"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +5.86%  "
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").Value = "'19.77"
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("D51").Value = "'1.86"
$ws.Range("E51").Value = "  -3.18%  "
